$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("M2").Value = 1.620350333333333
$ws.Range("N2").Value = 4.861051
$ws.Range("O2").Value = 0.0725197794467048
$ws.Range("P2").Value = 0.07251977944670479
$ws.Range("Q2").Value = 107.7159545728847
$ws.Range("R2").Value = 969.4435911559619
$ws.Range("S2").Value = 0.003127040638369827
$ws.Range("T2").Value = 0.003127040638369826
$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("O3").Value = 0.09715752300453648
$ws.Range("P3").Value = 0.09715752300453648
$ws.Range("Q3").Value = 144.3111853651147
$ws.Range("R3").Value = 1298.800668286032
$ws.Range("S3").Value = 0.004189415978323716
$ws.Range("T3").Value = 0.004189415978323716
$ws.Range("G4").Value = 66.47695399999999
$ws.Range("H4").Value = 199.430862
$ws.Range("I4").Value = 0.04311983106164722
$ws.Range("J4").Value = 0.04311983106164721
$ws.Range("M4").Value = 13.65106133333333
$ws.Range("N4").Value = 40.953184
$ws.Range("O4").Value = 0.6109616770777183
$ws.Range("P4").Value = 0.6109616770777183
$ws.Range("Q4").Value = 907.4809763071785
$ws.Range("R4").Value = 8167.328786764608
$ws.Range("S4").Value = 0.02634456430073188
$ws.Range("T4").Value = 0.02634456430073187
$ws.Range("G5").Value = 66.47695399999999
$ws.Range("H5").Value = 199.430862
$ws.Range("I5").Value = 0.04311983106164722
$ws.Range("J5").Value = 0.04311983106164721
$ws.Range("M5").Value = 0.6533493333333333
$ws.Range("N5").Value = 1.960048
$ws.Range("O5").Value = 0.02924105274043717
$ws.Range("P5").Value = 0.02924105274043717
$ws.Range("Q5").Value = 43.43267357793066
$ws.Range("R5").Value = 390.894062201376
$ws.Range("S5").Value = 0.001260869254232367
$ws.Range("T5").Value = 0.001260869254232367
$ws.Range("G6").Value = 66.47695399999999
$ws.Range("H6").Value = 199.430862
$ws.Range("I6").Value = 0.04311983106164722
$ws.Range("J6").Value = 0.04311983106164721
$ws.Range("M6").Value = 4.247957666666667
$ws.Range("N6").Value = 12.743873
$ws.Range("O6").Value = 0.1901199677306032
$ws.Range("P6").Value = 0.1901199677306032
$ws.Range("Q6").Value = 282.3912864009473
$ws.Range("R6").Value = 2541.521577608526
$ws.Range("S6").Value = 0.00819794088998943
$ws.Range("T6").Value = 0.008197940889989428
$ws.Range("I7").Value = 0.8830494168872806
$ws.Range("J7").Value = 0.8830494168872804
$ws.Range("M7").Value = 1.620350333333333
$ws.Range("N7").Value = 4.861051
$ws.Range("O7").Value = 0.0725197794467048
$ws.Range("P7").Value = 0.07251977944670479
$ws.Range("Q7").Value = 2205.911028247173
$ws.Range("R7").Value = 19853.19925422455
$ws.Range("S7").Value = 0.06403854895320688
$ws.Range("T7").Value = 0.06403854895320685
$ws.Range("I8").Value = 0.8830494168872806
$ws.Range("J8").Value = 0.8830494168872804
$ws.Range("O8").Value = 0.09715752300453648
$ws.Range("P8").Value = 0.09715752300453648
$ws.Range("S8").Value = 0.08579489403536848
$ws.Range("T8").Value = 0.08579489403536847
$ws.Range("I9").Value = 0.8830494168872806
$ws.Range("J9").Value = 0.8830494168872804
$ws.Range("M9").Value = 13.65106133333333
$ws.Range("N9").Value = 40.953184
$ws.Range("O9").Value = 0.6109616770777183
$ws.Range("P9").Value = 0.6109616770777183
$ws.Range("Q9").Value = 18584.26916883523
$ws.Range("R9").Value = 167258.4225195171
$ws.Range("S9").Value = 0.5395093526839542
$ws.Range("T9").Value = 0.5395093526839541
$ws.Range("I10").Value = 0.8830494168872806
$ws.Range("J10").Value = 0.8830494168872804
$ws.Range("M10").Value = 0.6533493333333333
$ws.Range("N10").Value = 1.960048
$ws.Range("O10").Value = 0.02924105274043717
$ws.Range("P10").Value = 0.02924105274043717
$ws.Range("Q10").Value = 889.4561071451039
$ws.Range("R10").Value = 8005.104964305936
$ws.Range("S10").Value = 0.02582129457161326
$ws.Range("T10").Value = 0.02582129457161326
$ws.Range("I11").Value = 0.8830494168872806
$ws.Range("J11").Value = 0.8830494168872804
$ws.Range("M11").Value = 4.247957666666667
$ws.Range("N11").Value = 12.743873
$ws.Range("O11").Value = 0.1901199677306032
$ws.Range("P11").Value = 0.1901199677306032
$ws.Range("Q11").Value = 5783.080653398079
$ws.Range("R11").Value = 52047.72588058271
$ws.Range("S11").Value = 0.1678853266431377
$ws.Range("T11").Value = 0.1678853266431377
$ws.Range("G12").Value = 44.831112
$ws.Range("H12").Value = 134.493336
$ws.Range("I12").Value = 0.02907940059566787
$ws.Range("J12").Value = 0.02907940059566786
$ws.Range("M12").Value = 1.620350333333333
$ws.Range("N12").Value = 4.861051
$ws.Range("O12").Value = 0.0725197794467048
$ws.Range("P12").Value = 0.07251977944670479
$ws.Range("Q12").Value = 72.642107272904
$ws.Range("R12").Value = 653.778965456136
$ws.Range("S12").Value = 0.002108831717640211
$ws.Range("T12").Value = 0.002108831717640209
$ws.Range("G13").Value = 44.831112
$ws.Range("H13").Value = 134.493336
$ws.Range("I13").Value = 0.02907940059566787
$ws.Range("J13").Value = 0.02907940059566786
$ws.Range("O13").Value = 0.09715752300453648
$ws.Range("P13").Value = 0.09715752300453648
$ws.Range("Q13").Value = 97.32141027334399
$ws.Range("R13").Value = 875.892692460096
$ws.Range("S13").Value = 0.002825282532331733
$ws.Range("T13").Value = 0.002825282532331732
$ws.Range("G14").Value = 44.831112
$ws.Range("H14").Value = 134.493336
$ws.Range("I14").Value = 0.02907940059566787
$ws.Range("J14").Value = 0.02907940059566786
$ws.Range("M14").Value = 13.65106133333333
$ws.Range("N14").Value = 40.953184
$ws.Range("O14").Value = 0.6109616770777183
$ws.Range("P14").Value = 0.6109616770777183
$ws.Range("Q14").Value = 611.992259553536
$ws.Range("R14").Value = 5507.930335981824
$ws.Range("S14").Value = 0.01776639935634404
$ws.Range("T14").Value = 0.01776639935634404
$ws.Range("G15").Value = 44.831112
$ws.Range("H15").Value = 134.493336
$ws.Range("I15").Value = 0.02907940059566787
$ws.Range("J15").Value = 0.02907940059566786
$ws.Range("M15").Value = 0.6533493333333333
$ws.Range("N15").Value = 1.960048
$ws.Range("O15").Value = 0.02924105274043717
$ws.Range("P15").Value = 0.02924105274043717
$ws.Range("Q15").Value = 29.290377137792
$ws.Range("R15").Value = 263.613394240128
$ws.Range("S15").Value = 0.0008503122864782243
$ws.Range("T15").Value = 0.0008503122864782241
$ws.Range("G16").Value = 44.831112
$ws.Range("H16").Value = 134.493336
$ws.Range("I16").Value = 0.02907940059566787
$ws.Range("J16").Value = 0.02907940059566786
$ws.Range("M16").Value = 4.247957666666667
$ws.Range("N16").Value = 12.743873
$ws.Range("O16").Value = 0.1901199677306032
$ws.Range("P16").Value = 0.1901199677306032
$ws.Range("Q16").Value = 190.440665925592
$ws.Range("R16").Value = 1713.965993330328
$ws.Range("S16").Value = 0.005528574702873659
$ws.Range("T16").Value = 0.005528574702873656
$ws.Range("G17").Value = 52.83062100000001
$ws.Range("H17").Value = 158.491863
$ws.Range("I17").Value = 0.0342682285413064
$ws.Range("J17").Value = 0.03426822854130639
$ws.Range("M17").Value = 1.620350333333333
$ws.Range("N17").Value = 4.861051
$ws.Range("O17").Value = 0.0725197794467048
$ws.Range("P17").Value = 0.07251977944670479
$ws.Range("Q17").Value = 85.60411434755702
$ws.Range("R17").Value = 770.4370291280131
$ws.Range("S17").Value = 0.002485124375844815
$ws.Range("T17").Value = 0.002485124375844814
$ws.Range("G18").Value = 52.83062100000001
$ws.Range("H18").Value = 158.491863
$ws.Range("I18").Value = 0.0342682285413064
$ws.Range("J18").Value = 0.03426822854130639
$ws.Range("O18").Value = 0.09715752300453648
$ws.Range("P18").Value = 0.09715752300453648
$ws.Range("Q18").Value = 114.687107054952
$ws.Range("R18").Value = 1032.183963494568
$ws.Range("S18").Value = 0.00332941620282669
$ws.Range("T18").Value = 0.003329416202826689
$ws.Range("G19").Value = 52.83062100000001
$ws.Range("H19").Value = 158.491863
$ws.Range("I19").Value = 0.0342682285413064
$ws.Range("J19").Value = 0.03426822854130639
$ws.Range("M19").Value = 13.65106133333333
$ws.Range("N19").Value = 40.953184
$ws.Range("O19").Value = 0.6109616770777183
$ws.Range("P19").Value = 0.6109616770777183
$ws.Range("Q19").Value = 721.1940475490881
$ws.Range("R19").Value = 6490.746427941793
$ws.Range("S19").Value = 0.02093657438007909
$ws.Range("T19").Value = 0.02093657438007909
$ws.Range("G20").Value = 52.83062100000001
$ws.Range("H20").Value = 158.491863
$ws.Range("I20").Value = 0.0342682285413064
$ws.Range("J20").Value = 0.03426822854130639
$ws.Range("M20").Value = 0.6533493333333333
$ws.Range("N20").Value = 1.960048
$ws.Range("O20").Value = 0.02924105274043717
$ws.Range("P20").Value = 0.02924105274043717
$ws.Range("Q20").Value = 34.516851009936
$ws.Range("R20").Value = 310.6516590894241
$ws.Range("S20").Value = 0.001002039078097695
$ws.Range("T20").Value = 0.001002039078097694
$ws.Range("G21").Value = 52.83062100000001
$ws.Range("H21").Value = 158.491863
$ws.Range("I21").Value = 0.0342682285413064
$ws.Range("J21").Value = 0.03426822854130639
$ws.Range("M21").Value = 4.247957666666667
$ws.Range("N21").Value = 12.743873
$ws.Range("O21").Value = 0.1901199677306032
$ws.Range("P21").Value = 0.1901199677306032
$ws.Range("Q21").Value = 224.4222415117111
$ws.Range("R21").Value = 2019.800173605399
$ws.Range("S21").Value = 0.006515074504458108
$ws.Range("T21").Value = 0.006515074504458106
$ws.Range("G22").Value = 16.16161433333333
$ws.Range("H22").Value = 48.484843
$ws.Range("I22").Value = 0.01048312291409786
$ws.Range("J22").Value = 0.01048312291409786
$ws.Range("M22").Value = 1.620350333333333
$ws.Range("N22").Value = 4.861051
$ws.Range("O22").Value = 0.0725197794467048
$ws.Range("P22").Value = 0.07251977944670479
$ws.Range("Q22").Value = 26.18747717222144
$ws.Range("R22").Value = 235.687294549993
$ws.Range("S22").Value = 0.0007602337616430745
$ws.Range("T22").Value = 0.0007602337616430742
$ws.Range("G23").Value = 16.16161433333333
$ws.Range("H23").Value = 48.484843
$ws.Range("I23").Value = 0.01048312291409786
$ws.Range("J23").Value = 0.01048312291409786
$ws.Range("O23").Value = 0.09715752300453648
$ws.Range("P23").Value = 0.09715752300453648
$ws.Range("Q23").Value = 35.08436505464977
$ws.Range("R23").Value = 315.759285491848
$ws.Range("S23").Value = 0.001018514255685847
$ws.Range("T23").Value = 0.001018514255685847
$ws.Range("G24").Value = 16.16161433333333
$ws.Range("H24").Value = 48.484843
$ws.Range("I24").Value = 0.01048312291409786
$ws.Range("J24").Value = 0.01048312291409786
$ws.Range("M24").Value = 13.65106133333333
$ws.Range("N24").Value = 40.953184
$ws.Range("O24").Value = 0.6109616770777183
$ws.Range("P24").Value = 0.6109616770777183
$ws.Range("Q24").Value = 220.6231885100124
$ws.Range("R24").Value = 1985.608696590112
$ws.Range("S24").Value = 0.006404786356609089
$ws.Range("T24").Value = 0.006404786356609088
$ws.Range("G25").Value = 16.16161433333333
$ws.Range("H25").Value = 48.484843
$ws.Range("I25").Value = 0.01048312291409786
$ws.Range("J25").Value = 0.01048312291409786
$ws.Range("M25").Value = 0.6533493333333333
$ws.Range("N25").Value = 1.960048
$ws.Range("O25").Value = 0.02924105274043717
$ws.Range("P25").Value = 0.02924105274043717
$ws.Range("Q25").Value = 10.55917995027378
$ws.Range("R25").Value = 95.032619552464
$ws.Range("S25").Value = 0.0003065375500156211
$ws.Range("T25").Value = 0.000306537550015621
$ws.Range("G26").Value = 16.16161433333333
$ws.Range("H26").Value = 48.484843
$ws.Range("I26").Value = 0.01048312291409786
$ws.Range("J26").Value = 0.01048312291409786
$ws.Range("M26").Value = 4.247957666666667
$ws.Range("N26").Value = 12.743873
$ws.Range("O26").Value = 0.1901199677306032
$ws.Range("P26").Value = 0.1901199677306032
$ws.Range("Q26").Value = 68.65385351299322
$ws.Range("R26").Value = 617.884681616939
$ws.Range("S26").Value = 0.001993050990144233
$ws.Range("T26").Value = 0.001993050990144232

$wb.Save()